$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 196.86957
$ws.Range("J33").Value = 496.66666
$ws.Range("L33").Value = 496.66666
$ws.Range("N33").Value = -954.66666

$ws.Range("H107").Value = 500
$ws.Range("I107").Value = 423.5
$ws.Range("J107").Value = 806
$ws.Range("K107").Value = 423.5
$ws.Range("L107").Value = 806
$ws.Range("M107").Value = 1496.5
$ws.Range("N107").Value = -4646

$ws.Range("H127").Value = 2452.889
$ws.Range("J127").Value = 2349.5
$ws.Range("L127").Value = 7048.5
$ws.Range("N127").Value = -16968.5

$ws.Range("H132").Value = 1947.3846
$ws.Range("I132").Value = 1984.6666
$ws.Range("K132").Value = 5953.9998
$ws.Range("M132").Value = -3423.9998

$ws.Range("H138").Value = 3786.0625
$ws.Range("J138").Value = 5800
$ws.Range("L138").Value = 17400
$ws.Range("N138").Value = -27680

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H32").Value = 3853327.8
$ws.Range("I32").Value = 3892553
$ws.Range("J32").Value = 3500300
$ws.Range("K32").Value = 3892553
$ws.Range("L32").Value = 3500300
$ws.Range("M32").Value = -3892266
$ws.Range("N32").Value = -3500874

$ws.Range("H45").Value = 3188.7144
$ws.Range("I45").Value = 2464.2666
$ws.Range("K45").Value = 2464.2666
$ws.Range("M45").Value = -2087.2666

$ws.Range("H61").Value = 10666.333
$ws.Range("I61").Value = 3499.5
$ws.Range("J61").Value = 25000
$ws.Range("K61").Value = 3499.5
$ws.Range("L61").Value = 25000
$ws.Range("M61").Value = -3287.5
$ws.Range("N61").Value = -25424

$ws.Range("H105").Value = 29999
$ws.Range("J105").Value = 29999
$ws.Range("L105").Value = 29999
$ws.Range("N105").Value = -36987

$ws.Range("H110").Value = 1820.6
$ws.Range("I110").Value = 2022.5
$ws.Range("K110").Value = 2022.5
$ws.Range("M110").Value = 22.5

$ws.Range("H132").Value = 1334.6364
$ws.Range("I132").Value = 1268.1
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 3804.3
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -1274.3
$ws.Range("N132").Value = -11060

$ws.Range("H136").Value = 10666.333
$ws.Range("I136").Value = 3499.5
$ws.Range("J136").Value = 25000
$ws.Range("K136").Value = 10498.5
$ws.Range("L136").Value = 75000
$ws.Range("M136").Value = -7948.5
$ws.Range("N136").Value = -80100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1816.4
$ws.Range("I105").Value = 1816.4
$ws.Range("K105").Value = 1816.4
$ws.Range("M105").Value = -69.40000000000009

$ws.Range("H107").Value = 2628.0557
$ws.Range("I107").Value = 2252.5
$ws.Range("K107").Value = 2252.5
$ws.Range("M107").Value = -332.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4069.2666
$ws.Range("I16").Value = 3044
$ws.Range("K16").Value = 3044
$ws.Range("M16").Value = -2757

$ws.Range("H58").Value = 2787.6924
$ws.Range("I58").Value = 1458.7142
$ws.Range("K58").Value = 1458.7142
$ws.Range("M58").Value = -1255.7142

$ws.Range("H105").Value = 2421.52
$ws.Range("I105").Value = 1767.1177
$ws.Range("J105").Value = 3812.125
$ws.Range("K105").Value = 1767.1177
$ws.Range("L105").Value = 3812.125
$ws.Range("M105").Value = -20.11770000000001
$ws.Range("N105").Value = -7306.125

$ws.Range("H107").Value = 314.2
$ws.Range("I107").Value = 257.33334
$ws.Range("J107").Value = 399.5
$ws.Range("K107").Value = 257.33334
$ws.Range("L107").Value = 399.5
$ws.Range("M107").Value = 1662.66666
$ws.Range("N107").Value = -4239.5

$ws.Range("H113").Value = 4069.2666
$ws.Range("I113").Value = 3044
$ws.Range("K113").Value = 3044
$ws.Range("M113").Value = -874

$ws.Range("H136").Value = 2787.6924
$ws.Range("I136").Value = 1458.7142
$ws.Range("K136").Value = 4376.142599999999
$ws.Range("M136").Value = -1826.142599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 10009548
$ws.Range("I4").Value = 10919053
$ws.Range("K4").Value = 32757159
$ws.Range("M4").Value = -32757047

$ws.Range("H68").Value = 2725.7954
$ws.Range("J68").Value = 2883.05
$ws.Range("L68").Value = 8649.150000000001
$ws.Range("N68").Value = -10271.15

$ws.Range("H71").Value = 2725.7954
$ws.Range("J71").Value = 2883.05
$ws.Range("L71").Value = 25947.45
$ws.Range("N71").Value = -34059.45

$ws.Range("H107").Value = 2119.7
$ws.Range("J107").Value = 2462.4285
$ws.Range("L107").Value = 7387.2855
$ws.Range("N107").Value = -11227.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3612.625
$ws.Range("I113").Value = 1760.4
$ws.Range("K113").Value = 1760.4
$ws.Range("M113").Value = 409.5999999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5214.143
$ws.Range("I7").Value = 4900
$ws.Range("K7").Value = 4900
$ws.Range("M7").Value = -4788

$ws.Range("H22").Value = 1670
$ws.Range("I22").Value = 550
$ws.Range("J22").Value = 1950
$ws.Range("K22").Value = 550
$ws.Range("L22").Value = 1950
$ws.Range("M22").Value = -255
$ws.Range("N22").Value = -2540

$ws.Range("H27").Value = 1670
$ws.Range("I27").Value = 550
$ws.Range("J27").Value = 1950
$ws.Range("K27").Value = 550
$ws.Range("L27").Value = 1950
$ws.Range("M27").Value = -443
$ws.Range("N27").Value = -2164

$ws.Range("H40").Value = 2707.9
$ws.Range("I40").Value = 2132.4
$ws.Range("J40").Value = 3283.4
$ws.Range("K40").Value = 2132.4
$ws.Range("L40").Value = 3283.4
$ws.Range("M40").Value = -1996.4
$ws.Range("N40").Value = -3555.4

$ws.Range("H61").Value = 1472.5
$ws.Range("I61").Value = 1296.8889
$ws.Range("K61").Value = 1296.8889
$ws.Range("M61").Value = -1094.8889

$ws.Range("H68").Value = 2598
$ws.Range("I68").Value = 2647.5
$ws.Range("J68").Value = 2499
$ws.Range("K68").Value = 2647.5
$ws.Range("L68").Value = 2499
$ws.Range("M68").Value = -1898.5
$ws.Range("N68").Value = -3997

$ws.Range("H71").Value = 2598
$ws.Range("I71").Value = 2647.5
$ws.Range("J71").Value = 2499
$ws.Range("K71").Value = 13237.5
$ws.Range("L71").Value = 12495
$ws.Range("M71").Value = -9493.5
$ws.Range("N71").Value = -19983

$ws.Range("H95").Value = 32248.75
$ws.Range("J95").Value = 32248.75
$ws.Range("L95").Value = 32248.75
$ws.Range("N95").Value = -37740.75

$ws.Range("H113").Value = 1472.5
$ws.Range("I113").Value = 1296.8889
$ws.Range("K113").Value = 1296.8889
$ws.Range("M113").Value = 873.1111000000001

$ws.Range("H126").Value = 5214.143
$ws.Range("I126").Value = 4900
$ws.Range("K126").Value = 14700
$ws.Range("M126").Value = -12230

$ws.Range("H135").Value = 89999
$ws.Range("J135").Value = 89999
$ws.Range("L135").Value = 89999
$ws.Range("N135").Value = -100139

$ws.Range("H136").Value = 1756.875
$ws.Range("J136").Value = 2000
$ws.Range("L136").Value = 6000
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 2249.5
$ws.Range("I14").Value = 1999
$ws.Range("K14").Value = 1999
$ws.Range("M14").Value = -1831

$ws.Range("H15").Value = 38079.332
$ws.Range("J15").Value = 38079.332
$ws.Range("L15").Value = 38079.332
$ws.Range("N15").Value = -38655.332

$ws.Range("H20").Value = 31024.25
$ws.Range("I20").Value = 4110
$ws.Range("J20").Value = 39995.668
$ws.Range("K20").Value = 4110
$ws.Range("L20").Value = 39995.668
$ws.Range("M20").Value = -3870
$ws.Range("N20").Value = -40475.668

$ws.Range("I39").Value = 10000
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 10000
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -9587
$ws.Range("N39").ClearContents()

$ws.Range("H107").Value = 2071.8572
$ws.Range("I107").Value = 1825.5
$ws.Range("J107").Value = 2400.3333
$ws.Range("K107").Value = 5476.5
$ws.Range("L107").Value = 7200.999899999999
$ws.Range("M107").Value = -3556.5
$ws.Range("N107").Value = -11040.9999

$ws.Range("H136").Value = 2318.7307
$ws.Range("I136").Value = 2339.28
$ws.Range("K136").Value = 7017.84
$ws.Range("M136").Value = -4467.84

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
